# Update cryptocurrency price (D) and volume-change (E) columns
# with freshly scraped values from coinranking.com
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D-column cells whose new text looks like a plain decimal number need to be
# pre-formatted as Text so Excel's COM auto-detection doesn't coerce the
# literal string (e.g. '1.000' or '2.690') into a numeric value.
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Write the new Price values (column D)
$ws.Range("D2").Value = '30.538.61'
$ws.Range("D3").Value = '1.875.72'
$ws.Range("D4").Value = '1.000'
$ws.Range("D5").Value = '236.01'
$ws.Range("D6").Value = '1.0000'
$ws.Range("D7").Value = '0.4874'
$ws.Range("D8").Value = '0.2898'
$ws.Range("D9").Value = '0.06661'
$ws.Range("D10").Value = '1.874.39'
$ws.Range("D11").Value = '16.54'
$ws.Range("D12").Value = '0.07237'
$ws.Range("D13").Value = '88.64'
$ws.Range("D14").Value = '4.994'
$ws.Range("D15").Value = '0.6496'
$ws.Range("D16").Value = '30.473.67'
$ws.Range("D17").Value = '0.000007849'
$ws.Range("D20").Value = '2.117.50'
$ws.Range("D21").Value = '1.001'
$ws.Range("D22").Value = '4.706'
$ws.Range("D23").Value = '192.49'
$ws.Range("D24").Value = '6.098'
$ws.Range("D25").Value = '9.319'
$ws.Range("D26").Value = '156.64'
$ws.Range("D27").Value = '18.37'
$ws.Range("D28").Value = '1.821'
$ws.Range("D29").Value = '1.406'
$ws.Range("D30").Value = '4.248'
$ws.Range("D31").Value = '0.09023'
$ws.Range("D32").Value = '3.914'
$ws.Range("D33").Value = '0.05108'
$ws.Range("D34").Value = '0.7206'
$ws.Range("D35").Value = '1.077'
$ws.Range("D36").Value = '2.690'
$ws.Range("D37").Value = '0.01814'
$ws.Range("D38").Value = '2.659'
$ws.Range("D39").Value = '0.9155'
$ws.Range("D41").Value = '0.4385'
$ws.Range("D42").Value = '104.56'
$ws.Range("D43").Value = '0.9944'
$ws.Range("D44").Value = '5.708'
$ws.Range("D45").Value = '0.1327'
$ws.Range("D46").Value = '7.372'
$ws.Range("D47").Value = '0.4022'
$ws.Range("D48").Value = '0.05815'
$ws.Range("D49").Value = '8.695'
$ws.Range("D50").Value = '1.401'
$ws.Range("D51").Value = '33.06'

# Write the new Volume(1h) values (column E)
$ws.Range("E2").Value = '  -0.19%  '
$ws.Range("E3").Value = '  -0.73%  '
$ws.Range("E4").Value = '  -0.32%  '
$ws.Range("E5").Value = '  -3.46%  '
$ws.Range("E6").Value = '  -0.31%  '
$ws.Range("E7").Value = '  -1.53%  '
$ws.Range("E8").Value = '  -1.63%  '
$ws.Range("E9").Value = '  -1.87%  '
$ws.Range("E10").Value = '  -0.79%  '
$ws.Range("E11").Value = '  -2.72%  '
$ws.Range("E12").Value = '  -0.81%  '
$ws.Range("E13").Value = '  -1.95%  '
$ws.Range("E14").Value = '  -0.90%  '
$ws.Range("E15").Value = '  -2.99%  '
$ws.Range("E16").Value = '  -0.49%  '
$ws.Range("E17").Value = '  -0.85%  '
$ws.Range("E18").Value = '  -0.11%  '
$ws.Range("E19").Value = '  -1.51%  '
$ws.Range("E20").Value = '  -1.11%  '
$ws.Range("E21").Value = '  -0.40%  '
$ws.Range("E22").Value = '  -2.81%  '
$ws.Range("E23").Value = '  +10.38%  '
$ws.Range("E24").Value = '  +1.01%  '
$ws.Range("E25").Value = '  +0.90%  '
$ws.Range("E26").Value = '  +0.91%  '
$ws.Range("E27").Value = '  -0.64%  '
$ws.Range("E28").Value = '  -5.10%  '
$ws.Range("E29").Value = '  +1.22%  '
$ws.Range("E30").Value = '  -1.67%  '
$ws.Range("E31").Value = '  +1.42%  '
$ws.Range("E32").Value = '  -2.36%  '
$ws.Range("E33").Value = '  -2.19%  '
$ws.Range("E34").Value = '  -2.27%  '
$ws.Range("E35").Value = '  -4.61%  '
$ws.Range("E36").Value = '  +0.59%  '
$ws.Range("E37").Value = '  -2.74%  '
$ws.Range("E38").Value = '  -1.40%  '
$ws.Range("E39").Value = '  -2.15%  '
$ws.Range("E40").Value = '  -5.20%  '
$ws.Range("E41").Value = '  +0.97%  '
$ws.Range("E42").Value = '  -1.09%  '
$ws.Range("E43").Value = '  -0.85%  '
$ws.Range("E44").Value = '  -1.56%  '
$ws.Range("E45").Value = '  -2.12%  '
$ws.Range("E46").Value = '  -3.44%  '
$ws.Range("E47").Value = '  +4.13%  '
$ws.Range("E48").Value = '  -0.26%  '
$ws.Range("E49").Value = '  +2.25%  '
$ws.Range("E50").Value = '  +1.82%  '
$ws.Range("E51").Value = '  -0.71%  '
